$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric need an explicit Text format first,
# otherwise Excel auto-converts strings like "1.000" or "10.04" into numbers
# and drops the significant trailing zeros / formatting that the source data uses.

$ws.Range("D2").Value = "30.281.18"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.931.45"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7519"
$ws.Range("E5").Value = "  +4.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.21"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.78"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3178"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07091"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7786"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08035"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "1.929.78"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.386"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.01"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.53"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "30.282.59"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.017"
$ws.Range("E18").Value = "  +4.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "251.87"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007935"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").Value = "2.194.25"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.686"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.529"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.62"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.09"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1301"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.191"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.540"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.410"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.133"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05216"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.316"
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7556"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01948"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.796"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.28"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.495"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4521"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.977"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8403"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.04"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.671"
$ws.Range("E47").Value = "  +3.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.62"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.88"
$ws.Range("E49").Value = "  +3.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1219"
$ws.Range("E50").Value = "  +7.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "960.78"
$ws.Range("E51").Value = "  +2.14%  "
